$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing existing rows 21-46 down to 22-47.
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with a new price record (mirrors the unchanged
# fields of the old row 21 that shifted to row 22, with updated values).
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("D21").Value = 44546
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103001
$ws.Range("J21").Value = "Cereza"
$ws.Range("K21").Value = "Lapins"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 500
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 9500
$ws.Range("P21").Value = 9250
$ws.Range("Q21").Value = "$/caja 8 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1156
$ws.Range("T21").Value = 8
